$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The two "백상지(A0)/(A1)" rows (671-672) were miscategorised under "7.기타".
# Save their B/C contents (column A will become "8.자재" at the new spot),
# then shift rows 673:685 up by two (671:683), and place the saved rows at
# the tail of the "7.기타" block (684:685) re-tagged as "8.자재".
$savedB671 = $ws.Range("B671").Value2
$savedC671 = $ws.Range("C671").Value2
$savedB672 = $ws.Range("B672").Value2
$savedC672 = $ws.Range("C672").Value2

for ($r = 671; $r -le 683; $r++) {
    $src = $r + 2
    $ws.Range("B$r").Value2 = $ws.Range("B$src").Value2
    $ws.Range("C$r").Value2 = $ws.Range("C$src").Value2
}

$ws.Range("A684").Value2 = $ws.Range("A686").Value2
$ws.Range("B684").Value2 = $savedB671
$ws.Range("C684").Value2 = $savedC671

$ws.Range("A685").Value2 = $ws.Range("A687").Value2
$ws.Range("B685").Value2 = $savedB672
$ws.Range("C685").Value2 = $savedC672

# Restore the cursor position captured the next time the workbook was saved.
$ws.Range("D11").Select()
